# Update "想去人数" (number of interested attendees) values for several
# exhibition rows. These updates apply identically to both the "展览"
# sheet and the "全部类型" sheet, which mirror the same data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3" = 2121
    "F4" = 1613
    "F6" = 1034
    "F7" = 500
    "F8" = 25
    "F9" = 5673
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
